$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the existing "From Query" table (Table4) so we can freely rewrite
#    its range, then rebuild the cell grid.
# ---------------------------------------------------------------------------
for ($i = $ws.ListObjects.Count; $i -ge 1; $i--) {
    $lo = $ws.ListObjects.Item($i)
    if ($lo.Name -eq "Table4") {
        $lo.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) Extend the "From Query" title merge F7:H7 -> F7:I7, and insert a new
#    "Class Type" column (with "Person" values) between "House Street" and
#    "Name", shifting the old Name/Age data one column to the right.
# ---------------------------------------------------------------------------
$ws.Range("F7:I7").Merge()
$ws.Range("F7").Value = "From Query"

$ws.Range("F8").Value = "House Street"
$ws.Range("G8").Value = "Class Type"
$ws.Range("H8").Value = "Name"
$ws.Range("I8").Value = "Age"

$ws.Range("F9").Value = "On Elm St."
$ws.Range("G9").Value = "Person"
$ws.Range("H9").Value = "John"
$ws.Range("I9").Value = 30

$ws.Range("F10").Value = "On 23rd St."
$ws.Range("G10").Value = "Person"
$ws.Range("H10").Value = "Luis"
$ws.Range("I10").Value = 21

$ws.Range("F11").Value = "On 5th Ave."
$ws.Range("G11").Value = "Person"
$ws.Range("H11").Value = "Henry"
$ws.Range("I11").Value = 45

# Copy the header/title formatting onto the newly-used I7 cell and make sure
# the new I8:I11 cells pick up the same plain style as the rest of the table.
$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("H9:H11").Copy()
$ws.Range("I9:I11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Add the new "From List" block: a title row, header row and three data
#    rows, mirroring the (now 4-column) "From Query" block above it.
# ---------------------------------------------------------------------------
$ws.Range("F15:I15").Merge()
$ws.Range("F15").Value = "From List"

$ws.Range("F16").Value = "House Street"
$ws.Range("G16").Value = "Class Type"
$ws.Range("H16").Value = "Name"
$ws.Range("I16").Value = "Age"

$ws.Range("F17").Value = "On Elm St."
$ws.Range("G17").Value = "Person"
$ws.Range("H17").Value = "John"
$ws.Range("I17").Value = 30

$ws.Range("F18").Value = "On 23rd St."
$ws.Range("G18").Value = "Person"
$ws.Range("H18").Value = "Luis"
$ws.Range("I18").Value = 21

$ws.Range("F19").Value = "On 5th Ave."
$ws.Range("G19").Value = "Person"
$ws.Range("H19").Value = "Henry"
$ws.Range("I19").Value = 45

# Copy formatting down from the "From Query" block onto the new "From List"
# block so the title/header/data rows look identical.
$ws.Range("F7:I7").Copy()
$ws.Range("F15:I15").PasteSpecial(-4122)
$ws.Range("F8:I8").Copy()
$ws.Range("F16:I16").PasteSpecial(-4122)
$ws.Range("F9:I11").Copy()
$ws.Range("F17:I19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Re-create the "From Query" table over its new F8:I11 extent, and add a
#    brand new table over the "From List" data at F16:I19.
# ---------------------------------------------------------------------------
$queryTable = $ws.ListObjects.Add(1, $ws.Range("F8:I11"), 0, 1)
$queryTable.Name = "Table4"
$queryTable.TableStyle = "TableStyleLight9"

$listTable = $ws.ListObjects.Add(1, $ws.Range("F16:I19"), 0, 1)
$listTable.Name = "Table5"
$listTable.TableStyle = "TableStyleLight9"

# ---------------------------------------------------------------------------
# 5) Column widths: the new "Class Type" column is a bit wider, and the new
#    trailing "Age" column is narrower.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 12.42
$ws.Columns.Item(9).ColumnWidth = 6.92

# ---------------------------------------------------------------------------
# 6) Update the "Titles" defined name to cover the extended F7:I7 range and
#    the new F15:I15 title range.
# ---------------------------------------------------------------------------
$wb.Names.Item("Titles").RefersTo = '=''Inserting Tables''!$A$1:$A$1,''Inserting Tables''!$C$1:$H$1,''Inserting Tables''!$A$7:$D$7,''Inserting Tables''!$F$7:$I$7,''Inserting Tables''!$F$15:$I$15'
